# Applies crypto price/volume updates per commit:
# "Updated cryptos list on Mon Mar 25 22:40:14 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.407.71"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.34%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.605.42"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.94%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.53"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.92"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.96%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.68%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.593.28"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.74%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.20%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.660"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.69%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.07"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +5.03%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.58%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.47%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.179.59"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.82%  "

# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.34"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.47%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.596.38"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.51%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.268.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.29%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.46"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.79%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.120"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.24%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.37"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.20"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +15.11%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +7.80%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.45"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.76%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.51%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.32%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +6.12%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.81%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.55"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +8.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.41%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "619.29"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.66%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.118"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.15%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.22"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.54%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.45%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.404"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.07"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.17%  "

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.07%  "

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.25%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.300.75"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +5.62%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.31%  "

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.67"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.29%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.17%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.19%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.13"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.55%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.71"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.88%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.08"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.04%  "

